# "Generate Report for Handoff" -- refresh the localization-status report
# with a new handoff (the old "38e3a257.../a7b6c30b..." source files are
# replaced by "5fc316cb.../ffff8c1cf609..."), and flip the per-language
# rows from "handed back" to "ready for handoff" / newly generated.

$wb = $excel.ActiveWorkbook

$oldFile1 = "38e3a257-1b23-4dd7-8e6a-6406b493b1cb.md"
$oldFile2 = "a7b6c30b-12e8-41b8-ac84-7067e67f547b.md"
$newFile1 = "5fc316cb-1163-4287-adb3-ae991cadfad9.md"
$newFile2 = "ffff8c1cf609-4f2a-4bde-928b-57d98bb3639e.md"

$zhXlf = "5fc316cb-1163-4287-adb3-ae991cadfad9.e8ba3ad45d87ef6fc7e4b43cb4485f17c8599e7a.zh-cn.xlf"
$deXlf = "5fc316cb-1163-4287-adb3-ae991cadfad9.e8ba3ad45d87ef6fc7e4b43cb4485f17c8599e7a.de-de.xlf"

$status = "Ready for handoff"
$genDateZh = "2016-08-25 15:06:42"
$genDateDe = "2016-08-25 15:06:47"
$overviewDate = "2016-08-25 15:06:47"
$epoch = "0001-01-01 00:00:00"

function Clear-TargetHandback($ws, $row) {
    # "Latest Target File" (I) / "Latest Handback File" (J) are wiped out
    # and lose their hyperlink + hyperlink styling as part of this handoff.
    $iAddr = "`$I`$" + $row
    foreach ($hl in @($ws.Hyperlinks)) {
        if ($hl.Range.Address() -eq $iAddr) {
            $hl.Delete()
        }
    }
    $ws.Range("I" + $row).ClearContents()
    $ws.Range("J" + $row).ClearContents()
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $newFile1
$ws.Range("E2").Value = $status
$ws.Range("F2").Value = $status
$ws.Range("G2").Value = $overviewDate

$ws.Range("A3").Value = $newFile2
$ws.Range("E3").Value = $status
$ws.Range("F3").Value = $status
$ws.Range("G3").Value = $overviewDate

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$B$2') {
        $hl.TextToDisplay = "e2e\" + $newFile1
    }
    if ($hl.Range.Address() -eq '$B$3') {
        $hl.TextToDisplay = "e2e\" + $newFile2
    }
}

$ws.Columns.Item(5).ColumnWidth = 17.08
$ws.Columns.Item(6).ColumnWidth = 17.08

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = $newFile1
$ws.Range("C2").Value = $status
$ws.Range("G2").Value = $zhXlf
$ws.Range("H2").Value = $genDateZh
$ws.Range("K2").Value = $epoch

$ws.Range("A3").Value = $newFile2
$ws.Range("C3").Value = $status
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = $zhXlf
$ws.Range("H3").Value = $genDateZh
$ws.Range("K3").Value = $epoch

Clear-TargetHandback $ws 2
Clear-TargetHandback $ws 3

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newFile1
    }
    if ($hl.Range.Address() -eq '$A$3') {
        $hl.TextToDisplay = $newFile2
    }
}

$ws.Columns.Item(3).ColumnWidth = 17.08
$ws.Columns.Item(9).ColumnWidth = 18.65
$ws.Columns.Item(10).ColumnWidth = 21.7

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = $newFile1
$ws.Range("C2").Value = $status
$ws.Range("G2").Value = $deXlf
$ws.Range("H2").Value = $genDateDe
$ws.Range("K2").Value = $epoch

$ws.Range("A3").Value = $newFile2
$ws.Range("C3").Value = $status
$ws.Range("F3").Value = "True"
$ws.Range("G3").Value = $deXlf
$ws.Range("H3").Value = $genDateDe
$ws.Range("K3").Value = $epoch

Clear-TargetHandback $ws 2
Clear-TargetHandback $ws 3

foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq '$A$2') {
        $hl.TextToDisplay = $newFile1
    }
    if ($hl.Range.Address() -eq '$A$3') {
        $hl.TextToDisplay = $newFile2
    }
}

$ws.Columns.Item(3).ColumnWidth = 17.08
$ws.Columns.Item(9).ColumnWidth = 18.65
$ws.Columns.Item(10).ColumnWidth = 21.7
